# Generate Report for Handoff
# Adds two new localization-status rows (5cd6925c... and ea2278a0...) to each
# of the three worksheets (Overview, zh-cn, de-de), growing every table from
# 2 data rows to 4 data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Insert a new row 3 (pushes the existing "7554db8a..." row down to row 4),
# inheriting formatting from row 2/3 the way Excel normally does.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "5cd6925c-69a9-41d8-a51a-6b758213e0be.md"
$ws.Range("B3").Value = "e2e\5cd6925c-69a9-41d8-a51a-6b758213e0be.md"
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-19 20:46:04"

# Append a new row 5 for "ea2278a0..."
$ws.Range("A5").Value = "ea2278a0-a226-4cf8-86bf-b748bd23d759.md"
$ws.Range("B5").Value = "e2e\ea2278a0-a226-4cf8-86bf-b748bd23d759.md"
$ws.Range("C5").Value = ".md"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-19 20:46:04"
$ws.Range("B5").Style = $ws.Range("B2").Style
$ws.Range("G5").Style = $ws.Range("G2").Style

# Rebuild hyperlinks in B2:B5 (row-insert does not re-anchor hyperlink ranges)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0732d96232f33a994e43e3c393203ee667d38041/e2e/6de80fa9-da02-4a91-962f-5c2fac1d6db7.md", "", "", "e2e\6de80fa9-da02-4a91-962f-5c2fac1d6db7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5cd6925c69a941d8a51a6b758213e0be/e2e/5cd6925c-69a9-41d8-a51a-6b758213e0be.md", "", "", "e2e\5cd6925c-69a9-41d8-a51a-6b758213e0be.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fbe3d4e7f1a88c4fc45b4559bdf163e73b0985a/e2e/7554db8a-11dc-4367-aeea-82a5511553a5.md", "", "", "e2e\7554db8a-11dc-4367-aeea-82a5511553a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea2278a0a2264cf886bfb748bd23d759/e2e/ea2278a0-a226-4cf8-86bf-b748bd23d759.md", "", "", "e2e\ea2278a0-a226-4cf8-86bf-b748bd23d759.md") | Out-Null

# Grow the "Overview" table to cover the two new rows
$loOverview = $ws.ListObjects.Item(1)
$loOverview.Resize($ws.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "5cd6925c-69a9-41d8-a51a-6b758213e0be.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "5cd6925c-69a9-41d8-a51a-6b758213e0be.6d9d03178d3b1fb07e02f1071b2845fd988f150c.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-19 20:45:57"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$ws.Range("A5:P5").Value = $ws.Range("A4:P4").Value
$ws.Range("A5").Value = "ea2278a0-a226-4cf8-86bf-b748bd23d759.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "False"
$ws.Range("G5").Value = "ea2278a0-a226-4cf8-86bf-b748bd23d759.0e49ae626e3e220ea15daa6246b04908bcf81958.zh-cn.xlf"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "False"
$ws.Range("P5").Value = ""

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0732d96232f33a994e43e3c393203ee667d38041/e2e/6de80fa9-da02-4a91-962f-5c2fac1d6db7.md", "", "", "6de80fa9-da02-4a91-962f-5c2fac1d6db7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d7d626864793516ed492afca0421424ce864e8ae/e2e/6de80fa9-da02-4a91-962f-5c2fac1d6db7.md", "", "", "6de80fa9-da02-4a91-962f-5c2fac1d6db7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5cd6925c69a941d8a51a6b758213e0be/e2e/5cd6925c-69a9-41d8-a51a-6b758213e0be.md", "", "", "5cd6925c-69a9-41d8-a51a-6b758213e0be.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fbe3d4e7f1a88c4fc45b4559bdf163e73b0985a/e2e/7554db8a-11dc-4367-aeea-82a5511553a5.md", "", "", "7554db8a-11dc-4367-aeea-82a5511553a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ea2278a0a2264cf886bfb748bd23d759/e2e/ea2278a0-a226-4cf8-86bf-b748bd23d759.md", "", "", "ea2278a0-a226-4cf8-86bf-b748bd23d759.md") | Out-Null

$loZhCn = $ws.ListObjects.Item(1)
$loZhCn.Resize($ws.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "5cd6925c-69a9-41d8-a51a-6b758213e0be.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "5cd6925c-69a9-41d8-a51a-6b758213e0be.6d9d03178d3b1fb07e02f1071b2845fd988f150c.de-de.xlf"
$ws.Range("H3").Value = "2016-08-19 20:46:04"
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "True"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "False"
$ws.Range("P3").Value = ""

$ws.Range("A5:P5").Value = $ws.Range("A4:P4").Value
$ws.Range("A5").Value = "ea2278a0-a226-4cf8-86bf-b748bd23d759.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "False"
$ws.Range("G5").Value = "ea2278a0-a226-4cf8-86bf-b748bd23d759.0e49ae626e3e220ea15daa6246b04908bcf81958.de-de.xlf"
$ws.Range("H5").Value = "2016-08-19 20:46:04"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "False"
$ws.Range("P5").Value = ""

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0732d96232f33a994e43e3c393203ee667d38041/e2e/6de80fa9-da02-4a91-962f-5c2fac1d6db7.md", "", "", "6de80fa9-da02-4a91-962f-5c2fac1d6db7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3168d2d9aec265e6f63384c92a600a55d9a3aea4/e2e/6de80fa9-da02-4a91-962f-5c2fac1d6db7.md", "", "", "6de80fa9-da02-4a91-962f-5c2fac1d6db7.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/5cd6925c69a941d8a51a6b758213e0be/e2e/5cd6925c-69a9-41d8-a51a-6b758213e0be.md", "", "", "5cd6925c-69a9-41d8-a51a-6b758213e0be.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2fbe3d4e7f1a88c4fc45b4559bdf163e73b0985a/e2e/7554db8a-11dc-4367-aeea-82a5511553a5.md", "", "", "7554db8a-11dc-4367-aeea-82a5511553a5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ea2278a0a2264cf886bfb748bd23d759/e2e/ea2278a0-a226-4cf8-86bf-b748bd23d759.md", "", "", "ea2278a0-a226-4cf8-86bf-b748bd23d759.md") | Out-Null

$loDeDe = $ws.ListObjects.Item(1)
$loDeDe.Resize($ws.Range("A1:P5"))
